# Add data for 2022-06-08
# Updates the 2022 year-to-date (column I) cumulative crime counts
# across the Citywide Totals, By Neighborhood, and individual
# neighborhood worksheets to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 2845
$ws.Range("I3").Value = 2935
$ws.Range("H4").Value = 1667
$ws.Range("I4").Value = 709
$ws.Range("I5").Value = 263
$ws.Range("I6").Value = 3360
$ws.Range("H7").Value = 25973
$ws.Range("I7").Value = 10112

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 94
$ws.Range("I4").Value = 40
$ws.Range("I6").Value = 70
$ws.Range("I8").Value = 641
$ws.Range("I11").Value = 163
$ws.Range("I15").Value = 132
$ws.Range("I19").Value = 277
$ws.Range("I20").Value = 259
$ws.Range("I23").Value = 89
$ws.Range("I24").Value = 27
$ws.Range("I25").Value = 47
$ws.Range("I29").Value = 676
$ws.Range("I33").Value = 469
$ws.Range("I36").Value = 136
$ws.Range("I37").Value = 326
$ws.Range("I42").Value = 347
$ws.Range("I47").Value = 74
$ws.Range("I50").Value = 41
$ws.Range("I52").Value = 212
$ws.Range("I53").Value = 117
$ws.Range("I54").Value = 229
$ws.Range("I55").Value = 108
$ws.Range("H63").Value = 199
$ws.Range("I63").Value = 41
$ws.Range("I64").Value = 91
$ws.Range("I65").Value = 223
$ws.Range("I67").Value = 398
$ws.Range("I71").Value = 27
$ws.Range("I73").Value = 85
$ws.Range("I78").Value = 138
$ws.Range("I79").Value = 254
$ws.Range("I83").Value = 203
$ws.Range("I85").Value = 465
$ws.Range("I86").Value = 57
$ws.Range("I89").Value = 110
$ws.Range("I91").Value = 123
$ws.Range("I93").Value = 58
$ws.Range("I95").Value = 164
$ws.Range("I98").Value = 60
$ws.Range("I99").Value = 182
$ws.Range("H101").Value = 25973
$ws.Range("I101").Value = 10112

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 120
$ws.Range("I7").Value = 465

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 212

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 163

# Sheet 7: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 201
$ws.Range("I6").Value = 204
$ws.Range("I7").Value = 641

# Sheet 8: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 23
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 117

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 110

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 105
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 326

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 46
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 182

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 91
$ws.Range("I3").Value = 136
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 141
$ws.Range("I7").Value = 398

# Sheet 19: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 59
$ws.Range("I7").Value = 223

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 203

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 164

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 109
$ws.Range("I5").Value = 13
$ws.Range("I7").Value = 469

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 55
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 229

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 209
$ws.Range("I3").Value = 234
$ws.Range("I6").Value = 186
$ws.Range("I7").Value = 676

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 74
$ws.Range("I7").Value = 277

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I5").Value = 1
$ws.Range("I7").Value = 70

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 93
$ws.Range("I3").Value = 115
$ws.Range("I7").Value = 347

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 34
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 138

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 108

# Sheet 37: Dunning
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 27

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 89

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 43
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 123

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 254

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 91

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 72
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 259

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 136

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 58

# Sheet 52: East Side
$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 47

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 12
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 74

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 43
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 132

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 60

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 41

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 85

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 94

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 35
$ws.Range("I7").Value = 57

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 27

# Sheet 90: Archer Heights
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 40
